# Auto-generated edit script: refreshes currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# sheets to match the latest market-board snapshot pulled by the scheduled
# runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 37569.54
$ws.Range("I6").Value = 48388.7
$ws.Range("K6").Value = 145166.1
$ws.Range("M6").Value = -145054.1
$ws.Range("H12").Value = 2100
$ws.Range("I12").Value = 233.66667
$ws.Range("K12").Value = 233.66667
$ws.Range("M12").Value = -63.66667000000001
$ws.Range("H17").Value = 3055.0732
$ws.Range("J17").Value = 3680.5454
$ws.Range("L17").Value = 11041.6362
$ws.Range("N17").Value = -11377.6362
$ws.Range("H28").Value = 55724.633
$ws.Range("I28").Value = 64610.125
$ws.Range("K28").Value = 64610.125
$ws.Range("M28").Value = -64125.125
$ws.Range("H29").Value = 4748.1
$ws.Range("J29").Value = 7201.6
$ws.Range("L29").Value = 21604.8
$ws.Range("N29").Value = -22166.8
$ws.Range("H98").Value = 900.9091
$ws.Range("I98").Value = 900.9091
$ws.Range("K98").Value = 900.9091
$ws.Range("M98").Value = 597.0909
$ws.Range("H112").Value = 1551.1333
$ws.Range("I112").Value = 1157.8
$ws.Range("J112").Value = 1747.8
$ws.Range("K112").Value = 3473.4
$ws.Range("L112").Value = 5243.4
$ws.Range("M112").Value = -2365.4
$ws.Range("N112").Value = -7459.4
$ws.Range("H122").Value = 900.9091
$ws.Range("I122").Value = 900.9091
$ws.Range("K122").Value = 2702.7273
$ws.Range("M122").Value = -252.7273
$ws.Range("H137").Value = 2450.96
$ws.Range("I137").Value = 1836.2307
$ws.Range("J137").Value = 3116.9167
$ws.Range("K137").Value = 5508.6921
$ws.Range("L137").Value = 9350.750100000001
$ws.Range("M137").Value = -2958.6921
$ws.Range("N137").Value = -14450.7501
$ws.Range("H138").Value = 2894.6511
$ws.Range("I138").Value = 2063.0435
$ws.Range("J138").Value = 3851
$ws.Range("K138").Value = 6189.130500000001
$ws.Range("L138").Value = 11553
$ws.Range("M138").Value = -1049.130500000001
$ws.Range("N138").Value = -21833

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4214.4863
$ws.Range("I32").Value = 3590.5625
$ws.Range("J32").Value = 8207.6
$ws.Range("K32").Value = 3590.5625
$ws.Range("L32").Value = 8207.6
$ws.Range("M32").Value = -3303.5625
$ws.Range("N32").Value = -8781.6
$ws.Range("H45").Value = 45456984
$ws.Range("I45").Value = 66667950
$ws.Range("K45").Value = 66667950
$ws.Range("M45").Value = -66667573
$ws.Range("H46").Value = 11369.5625
$ws.Range("J46").Value = 4681.769
$ws.Range("L46").Value = 4681.769
$ws.Range("N46").Value = -5319.769
$ws.Range("H63").Value = 5844.4443
$ws.Range("I63").Value = 2650
$ws.Range("J63").Value = 8400
$ws.Range("K63").Value = 2650
$ws.Range("L63").Value = 8400
$ws.Range("M63").Value = -1964
$ws.Range("N63").Value = -9772
$ws.Range("H66").Value = 5844.4443
$ws.Range("I66").Value = 2650
$ws.Range("J66").Value = 8400
$ws.Range("K66").Value = 13250
$ws.Range("L66").Value = 42000
$ws.Range("M66").Value = -9818
$ws.Range("N66").Value = -48864
$ws.Range("H110").Value = 7622.9
$ws.Range("J110").Value = 10406.5
$ws.Range("L110").Value = 10406.5
$ws.Range("N110").Value = -14496.5
$ws.Range("H122").Value = 3055.8928
$ws.Range("I122").Value = 2187.6843
$ws.Range("K122").Value = 6563.0529
$ws.Range("M122").Value = -4113.0529

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1646
$ws.Range("I20").Value = 1061.2
$ws.Range("J20").Value = 2377
$ws.Range("K20").Value = 1061.2
$ws.Range("L20").Value = 2377
$ws.Range("M20").Value = -814.2
$ws.Range("N20").Value = -2871
$ws.Range("H99").Value = 2290
$ws.Range("I99").Value = 2703.8
$ws.Range("K99").Value = 2703.8
$ws.Range("M99").Value = -1205.8
$ws.Range("H105").Value = 53079.8
$ws.Range("I105").Value = 200000
$ws.Range("K105").Value = 200000
$ws.Range("M105").Value = -198253
$ws.Range("H107").Value = 3476.0908
$ws.Range("I107").Value = 2977
$ws.Range("J107").Value = 4807
$ws.Range("K107").Value = 2977
$ws.Range("L107").Value = 4807
$ws.Range("M107").Value = -1057
$ws.Range("N107").Value = -8647

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 582
$ws.Range("I16").Value = 633.3333
$ws.Range("J16").Value = 489.6
$ws.Range("K16").Value = 633.3333
$ws.Range("L16").Value = 489.6
$ws.Range("M16").Value = -346.3333
$ws.Range("N16").Value = -1063.6
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H86").Value = 14373.375
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877
$ws.Range("H89").Value = 14373.375
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384
$ws.Range("H107").Value = 1695.5
$ws.Range("I107").Value = 3125
$ws.Range("J107").Value = 1338.125
$ws.Range("K107").Value = 3125
$ws.Range("L107").Value = 1338.125
$ws.Range("M107").Value = -1205
$ws.Range("N107").Value = -5178.125
$ws.Range("H113").Value = 582
$ws.Range("I113").Value = 633.3333
$ws.Range("J113").Value = 489.6
$ws.Range("K113").Value = 633.3333
$ws.Range("L113").Value = 489.6
$ws.Range("M113").Value = 1536.6667
$ws.Range("N113").Value = -4829.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 41774.668
$ws.Range("I7").Value = 107.8
$ws.Range("J7").Value = 71536.71000000001
$ws.Range("K7").Value = 323.4
$ws.Range("L7").Value = 214610.13
$ws.Range("M7").Value = -211.4
$ws.Range("N7").Value = -214834.13
$ws.Range("H17").Value = 461
$ws.Range("I17").Value = 250
$ws.Range("J17").Value = 503.2
$ws.Range("K17").Value = 750
$ws.Range("L17").Value = 1509.6
$ws.Range("M17").Value = -581
$ws.Range("N17").Value = -1847.6
$ws.Range("H121").Value = 2800.923
$ws.Range("I121").Value = 1890
$ws.Range("J121").Value = 3370.25
$ws.Range("K121").Value = 5670
$ws.Range("L121").Value = 10110.75
$ws.Range("M121").Value = -4360
$ws.Range("N121").Value = -12730.75
$ws.Range("H122").Value = 1795.7222
$ws.Range("I122").Value = 292
$ws.Range("J122").Value = 2096.4666
$ws.Range("K122").Value = 2628
$ws.Range("L122").Value = 18868.1994
$ws.Range("M122").Value = -178
$ws.Range("N122").Value = -23768.1994
$ws.Range("H126").Value = 956
$ws.Range("I126").Value = 956
$ws.Range("K126").Value = 2868
$ws.Range("M126").Value = 2072
$ws.Range("H138").Value = 5692.6665
$ws.Range("I138").Value = 2652.625
$ws.Range("K138").Value = 7957.875
$ws.Range("M138").Value = -2817.875
$ws.Range("H140").Value = 2035.7646
$ws.Range("I140").Value = 1408.1428
$ws.Range("J140").Value = 4964.6665
$ws.Range("K140").Value = 4224.428400000001
$ws.Range("L140").Value = 14893.9995
$ws.Range("M140").Value = 955.5715999999993
$ws.Range("N140").Value = -25253.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9799.799999999999
$ws.Range("I70").Value = 9499.5
$ws.Range("K70").Value = 9499.5
$ws.Range("M70").Value = -9229.5
$ws.Range("H73").Value = 9799.799999999999
$ws.Range("I73").Value = 9499.5
$ws.Range("K73").Value = 9499.5
$ws.Range("M73").Value = -8563.5
$ws.Range("H107").Value = 2447.4285
$ws.Range("I107").Value = 1362.25
$ws.Range("J107").Value = 3894.3333
$ws.Range("K107").Value = 1362.25
$ws.Range("L107").Value = 3894.3333
$ws.Range("M107").Value = 557.75
$ws.Range("N107").Value = -7734.3333
$ws.Range("H123").Value = 27468.857
$ws.Range("J123").Value = 27468.857
$ws.Range("L123").Value = 27468.857
$ws.Range("N123").Value = -32368.857

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2381622.5
$ws.Range("I55").Value = 3125520.2
$ws.Range("J55").Value = 1150
$ws.Range("K55").Value = 3125520.2
$ws.Range("L55").Value = 1150
$ws.Range("M55").Value = -3125347.2
$ws.Range("N55").Value = -1496
$ws.Range("H100").Value = 3533.2068
$ws.Range("I100").Value = 2477.1304
$ws.Range("J100").Value = 7581.5
$ws.Range("K100").Value = 2477.1304
$ws.Range("L100").Value = 7581.5
$ws.Range("M100").Value = -1936.1304
$ws.Range("N100").Value = -8663.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2959.8096
$ws.Range("I81").Value = 1989.6154
$ws.Range("K81").Value = 3979.2308
$ws.Range("M81").Value = -2918.2308
$ws.Range("H84").Value = 2959.8096
$ws.Range("I84").Value = 1989.6154
$ws.Range("K84").Value = 19896.154
$ws.Range("M84").Value = -14592.154
$ws.Range("H107").Value = 1230.5714
$ws.Range("I107").Value = 1017.25
$ws.Range("K107").Value = 3051.75
$ws.Range("M107").Value = -1131.75
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 467
$ws.Range("I113").Value = 390.85715
$ws.Range("K113").Value = 1172.57145
$ws.Range("M113").Value = 997.4285500000001
$ws.Range("H122").Value = 5194.8335
$ws.Range("I122").Value = 3939.3635
$ws.Range("K122").Value = 11818.0905
$ws.Range("M122").Value = -9368.0905
$ws.Range("H132").Value = 6796.2666
$ws.Range("I132").Value = 6810.75
$ws.Range("K132").Value = 20432.25
$ws.Range("M132").Value = -17902.25

